$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 768.8
$ws.Cells.Item(15, 9).Value = 768.8
$ws.Cells.Item(15, 11).Value = 2306.4
$ws.Cells.Item(15, 13).Value = -2137.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 463.18182
$ws.Cells.Item(28, 9).Value = 480.7143
$ws.Cells.Item(28, 10).Value = 432.5
$ws.Cells.Item(28, 11).Value = 480.7143
$ws.Cells.Item(28, 12).Value = 432.5
$ws.Cells.Item(28, 13).Value = 4.28570000000002
$ws.Cells.Item(28, 14).Value = -1402.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(75, 8).Value = 40314
$ws.Cells.Item(75, 10).Value = 40314
$ws.Cells.Item(75, 12).Value = 40314
$ws.Cells.Item(75, 14).Value = -42186

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(78, 8).Value = 40314
$ws.Cells.Item(78, 10).Value = 40314
$ws.Cells.Item(78, 12).Value = 120942
$ws.Cells.Item(78, 14).Value = -130302

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 693.3461
$ws.Cells.Item(107, 9).Value = 534.25
$ws.Cells.Item(107, 10).Value = 947.9
$ws.Cells.Item(107, 11).Value = 534.25
$ws.Cells.Item(107, 12).Value = 947.9
$ws.Cells.Item(107, 13).Value = 1385.75
$ws.Cells.Item(107, 14).Value = -4787.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 2946.5386
$ws.Cells.Item(113, 9).Value = 2595.5557
$ws.Cells.Item(113, 10).Value = 3736.25
$ws.Cells.Item(113, 11).Value = 2595.5557
$ws.Cells.Item(113, 12).Value = 3736.25
$ws.Cells.Item(113, 13).Value = 658.4443000000001
$ws.Cells.Item(113, 14).Value = -10244.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 759.125
$ws.Cells.Item(129, 9).Value = 329
$ws.Cells.Item(129, 10).Value = 902.5
$ws.Cells.Item(129, 11).Value = 987
$ws.Cells.Item(129, 12).Value = 2707.5
$ws.Cells.Item(129, 13).Value = 4013
$ws.Cells.Item(129, 14).Value = -12707.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 5917.4707
$ws.Cells.Item(137, 9).Value = 5842.4287
$ws.Cells.Item(137, 10).Value = 6267.6665
$ws.Cells.Item(137, 11).Value = 17527.2861
$ws.Cells.Item(137, 12).Value = 18802.9995
$ws.Cells.Item(137, 13).Value = -14977.2861
$ws.Cells.Item(137, 14).Value = -23902.9995

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 5637.778
$ws.Cells.Item(141, 9).Value = 6105
$ws.Cells.Item(141, 10).Value = 1900
$ws.Cells.Item(141, 11).Value = 18315
$ws.Cells.Item(141, 12).Value = 5700
$ws.Cells.Item(141, 13).Value = -13135
$ws.Cells.Item(141, 14).Value = -16060

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 828.5714
$ws.Cells.Item(2, 9).Value = 734.0833
$ws.Cells.Item(2, 10).Value = 1395.5
$ws.Cells.Item(2, 11).Value = 734.0833
$ws.Cells.Item(2, 12).Value = 1395.5
$ws.Cells.Item(2, 13).Value = -621.0833
$ws.Cells.Item(2, 14).Value = -1621.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1407.5
$ws.Cells.Item(45, 9).Value = 1371.8182
$ws.Cells.Item(45, 10).Value = 1800
$ws.Cells.Item(45, 11).Value = 1371.8182
$ws.Cells.Item(45, 12).Value = 1800
$ws.Cells.Item(45, 13).Value = -994.8181999999999
$ws.Cells.Item(45, 14).Value = -2554

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 692.1070999999999
$ws.Cells.Item(97, 9).Value = 728.7083
$ws.Cells.Item(97, 11).Value = 728.7083
$ws.Cells.Item(97, 13).Value = -232.7083

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 828.5714
$ws.Cells.Item(116, 9).Value = 734.0833
$ws.Cells.Item(116, 10).Value = 1395.5
$ws.Cells.Item(116, 11).Value = 734.0833
$ws.Cells.Item(116, 12).Value = 1395.5
$ws.Cells.Item(116, 13).Value = 1559.9167
$ws.Cells.Item(116, 14).Value = -5983.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 828.5714
$ws.Cells.Item(3, 9).Value = 734.0833
$ws.Cells.Item(3, 10).Value = 1395.5
$ws.Cells.Item(3, 11).Value = 734.0833
$ws.Cells.Item(3, 12).Value = 1395.5
$ws.Cells.Item(3, 13).Value = -620.0833
$ws.Cells.Item(3, 14).Value = -1623.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 1351.5714
$ws.Cells.Item(99, 9).Value = 1196.5
$ws.Cells.Item(99, 10).Value = 1603.5625
$ws.Cells.Item(99, 11).Value = 1196.5
$ws.Cells.Item(99, 12).Value = 1603.5625
$ws.Cells.Item(99, 13).Value = 301.5
$ws.Cells.Item(99, 14).Value = -4599.5625

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 1351.5714
$ws.Cells.Item(126, 9).Value = 1196.5
$ws.Cells.Item(126, 10).Value = 1603.5625
$ws.Cells.Item(126, 11).Value = 3589.5
$ws.Cells.Item(126, 12).Value = 4810.6875
$ws.Cells.Item(126, 13).Value = -1119.5
$ws.Cells.Item(126, 14).Value = -9750.6875

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 108.75
$ws.Cells.Item(23, 9).Value = 95
$ws.Cells.Item(23, 10).Value = 110.71429
$ws.Cells.Item(23, 11).Value = 285
$ws.Cells.Item(23, 12).Value = 332.14287
$ws.Cells.Item(23, 13).Value = -50
$ws.Cells.Item(23, 14).Value = -802.14287

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 1529.6
$ws.Cells.Item(39, 10).Value = 1603.1428
$ws.Cells.Item(39, 12).Value = 4809.428400000001
$ws.Cells.Item(39, 14).Value = -5397.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(81, 8).Value = 2098.3333
$ws.Cells.Item(81, 10).Value = 2412.1428
$ws.Cells.Item(81, 12).Value = 7236.428400000001
$ws.Cells.Item(81, 14).Value = -9482.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(84, 8).Value = 2098.3333
$ws.Cells.Item(84, 10).Value = 2412.1428
$ws.Cells.Item(84, 12).Value = 21709.2852
$ws.Cells.Item(84, 14).Value = -32941.2852

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(99, 8).Value = 2490.2144
$ws.Cells.Item(99, 9).Value = 1913.8889
$ws.Cells.Item(99, 11).Value = 5741.6667
$ws.Cells.Item(99, 13).Value = -3495.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 849.322
$ws.Cells.Item(131, 9).Value = 475.875
$ws.Cells.Item(131, 10).Value = 907.902
$ws.Cells.Item(131, 11).Value = 1427.625
$ws.Cells.Item(131, 12).Value = 2723.706
$ws.Cells.Item(131, 13).Value = 3612.375
$ws.Cells.Item(131, 14).Value = -12803.706

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value = 18448.5
$ws.Cells.Item(57, 9).Value = 10000
$ws.Cells.Item(57, 10).Value = 19387.223
$ws.Cells.Item(57, 11).Value = 10000
$ws.Cells.Item(57, 12).Value = 19387.223
$ws.Cells.Item(57, 13).Value = -9180
$ws.Cells.Item(57, 14).Value = -21027.223

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 676.88464
$ws.Cells.Item(113, 9).Value = 573.86957
$ws.Cells.Item(113, 10).Value = 1466.6666
$ws.Cells.Item(113, 11).Value = 573.86957
$ws.Cells.Item(113, 12).Value = 1466.6666
$ws.Cells.Item(113, 13).Value = 1596.13043
$ws.Cells.Item(113, 14).Value = -5806.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2614.353
$ws.Cells.Item(122, 9).Value = 1060.75
$ws.Cells.Item(122, 10).Value = 3995.3333
$ws.Cells.Item(122, 11).Value = 3182.25
$ws.Cells.Item(122, 12).Value = 11985.9999
$ws.Cells.Item(122, 13).Value = -732.25
$ws.Cells.Item(122, 14).Value = -16885.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2819.889
$ws.Cells.Item(40, 9).Value = 2796
$ws.Cells.Item(40, 11).Value = 2796
$ws.Cells.Item(40, 13).Value = -2660

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(94, 8).Value = 9715
$ws.Cells.Item(94, 10).Value = 9715
$ws.Cells.Item(94, 12).Value = 9715
$ws.Cells.Item(94, 14).Value = -11067

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(76, 8).Value = 19999.25
$ws.Cells.Item(76, 10).Value = 19999.25
$ws.Cells.Item(76, 12).Value = 19999.25
$ws.Cells.Item(76, 14).Value = -20629.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(79, 8).Value = 19999.25
$ws.Cells.Item(79, 10).Value = 19999.25
$ws.Cells.Item(79, 12).Value = 19999.25
$ws.Cells.Item(79, 14).Value = -20629.25
